# Weekly update: a new Mango price record for "Feria Lagunitas de Puerto
# Montt" is prepended as the new row 100 (date 2022-04-18 / serial 44669),
# pushing the existing rows 100-190 down to 101-191.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 100; everything below shifts down one row.
$ws.Rows.Item(100).Insert()

# Populate the new row with the latest weekly price record.
$ws.Range("A100").Value = 4
$ws.Range("B100").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C100").Value = "Los Lagos"
$ws.Range("D100").Value2 = 44669
$ws.Range("E100").Value = 10
$ws.Range("F100").Value = "Fruta"
$ws.Range("G100").Value = 100108
$ws.Range("H100").Value = "Tropicales y subtropicales"
$ws.Range("I100").Value = 100108002
$ws.Range("J100").Value = "Mango"
$ws.Range("K100").Value = "Sin especificar"
$ws.Range("L100").Value = "Primera"
$ws.Range("M100").Value = 60
$ws.Range("N100").Value = 7500
$ws.Range("O100").Value = 8000
$ws.Range("P100").Value = 7750
$ws.Range("Q100").Value = "$/bandeja 4 kilos"
$ws.Range("R100").Value = "Perú"
$ws.Range("S100").Value = 1938
$ws.Range("T100").Value = 4
